# Auto-generated Excel COM-interop script applying the Sagittarius_Profits diff
# (workbook tab names ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR correspond to the
# per-sheet 'Leve' profit tables referenced by the diff, concatenated in the
# commit's canonical XML in workbook order.)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2754.1765
$ws.Range("I9").Value = 3707.2727
$ws.Range("K9").Value = 3707.2727
$ws.Range("M9").Value = -3538.2727

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H52").Value = 275
$ws.Range("I52").Value = 250
$ws.Range("K52").Value = 750
$ws.Range("M52").Value = -590

$ws.Range("H70").Value = 45773.12
$ws.Range("J70").Value = 64524
$ws.Range("L70").Value = 193572
$ws.Range("N70").Value = -194112

$ws.Range("H73").Value = 45773.12
$ws.Range("J73").Value = 64524
$ws.Range("L73").Value = 193572
$ws.Range("N73").Value = -195444

$ws.Range("H106").Value = 2237.7368
$ws.Range("I106").Value = 2323.2222
$ws.Range("J106").Value = 699
$ws.Range("K106").Value = 2323.2222
$ws.Range("L106").Value = 699
$ws.Range("M106").Value = -1692.2222
$ws.Range("N106").Value = -1961

$ws.Range("H132").Value = 1506.3334
$ws.Range("I132").Value = 1430.7
$ws.Range("K132").Value = 4292.1
$ws.Range("M132").Value = -1762.1

$ws.Range("H138").Value = 3754.86
$ws.Range("I138").Value = 3507.7
$ws.Range("J138").Value = 3816.65
$ws.Range("K138").Value = 10523.1
$ws.Range("L138").Value = 11449.95
$ws.Range("M138").Value = -5383.099999999999
$ws.Range("N138").Value = -21729.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 1099
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

$ws.Range("H39").Value = 1070.6666
$ws.Range("I39").Value = 1070.6666
$ws.Range("K39").Value = 1070.6666
$ws.Range("M39").Value = -550.6666

$ws.Range("H74").Value = 2183.1738
$ws.Range("I74").Value = 2029.7561
$ws.Range("K74").Value = 2029.7561
$ws.Range("M74").Value = -1155.7561

$ws.Range("H77").Value = 2183.1738
$ws.Range("I77").Value = 2029.7561
$ws.Range("K77").Value = 10148.7805
$ws.Range("M77").Value = -5780.780500000001

$ws.Range("H102").Value = 1647.6666
$ws.Range("I102").Value = 1709.2858
$ws.Range("J102").Value = 785
$ws.Range("K102").Value = 1709.2858
$ws.Range("L102").Value = 785
$ws.Range("M102").Value = -87.28580000000011
$ws.Range("N102").Value = -4029

$ws.Range("H110").Value = 100
$ws.Range("I110").Value = 100
$ws.Range("K110").Value = 100
$ws.Range("M110").Value = 1945

$ws.Range("H132").Value = 2019.1786
$ws.Range("I132").Value = 1924.5
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 5773.5
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -3243.5
$ws.Range("N132").Value = -14810

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2758.7273
$ws.Range("I20").Value = 2864.6
$ws.Range("J20").Value = 1700
$ws.Range("K20").Value = 2864.6
$ws.Range("L20").Value = 1700
$ws.Range("M20").Value = -2617.6
$ws.Range("N20").Value = -2194

$ws.Range("H86").Value = 2949.7
$ws.Range("I86").Value = 2220.6
$ws.Range("K86").Value = 2220.6
$ws.Range("M86").Value = -1097.6

$ws.Range("H89").Value = 2949.7
$ws.Range("I89").Value = 2220.6
$ws.Range("K89").Value = 11103
$ws.Range("M89").Value = -5487

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2381.923
$ws.Range("I16").Value = 2662.125
$ws.Range("J16").Value = 1933.6
$ws.Range("K16").Value = 2662.125
$ws.Range("L16").Value = 1933.6
$ws.Range("M16").Value = -2375.125
$ws.Range("N16").Value = -2507.6

$ws.Range("H22").Value = 9364.179
$ws.Range("I22").Value = 325.05554
$ws.Range("J22").Value = 25634.6
$ws.Range("K22").Value = 325.05554
$ws.Range("L22").Value = 25634.6
$ws.Range("M22").Value = 24.94445999999999
$ws.Range("N22").Value = -26334.6

$ws.Range("H62").Value = 3597.4
$ws.Range("I62").Value = 2994.5
$ws.Range("K62").Value = 2994.5
$ws.Range("M62").Value = -2370.5

$ws.Range("H65").Value = 3597.4
$ws.Range("I65").Value = 2994.5
$ws.Range("K65").Value = 14972.5
$ws.Range("M65").Value = -11852.5

$ws.Range("H113").Value = 2381.923
$ws.Range("I113").Value = 2662.125
$ws.Range("J113").Value = 1933.6
$ws.Range("K113").Value = 2662.125
$ws.Range("L113").Value = 1933.6
$ws.Range("M113").Value = -492.125
$ws.Range("N113").Value = -6273.6

$ws.Range("H134").Value = 3385.0715
$ws.Range("I134").Value = 3299.25
$ws.Range("K134").Value = 9897.75
$ws.Range("M134").Value = -7362.75

$ws.Range("H141").Value = 296353.25
$ws.Range("J141").Value = 296353.25
$ws.Range("L141").Value = 296353.25
$ws.Range("N141").Value = -306713.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 992.2
$ws.Range("J29").Value = 1499.3334
$ws.Range("L29").Value = 4498.0002
$ws.Range("N29").Value = -5052.0002

$ws.Range("H31").Value = 2300
$ws.Range("J31").Value = 2300
$ws.Range("L31").Value = 6900
$ws.Range("N31").Value = -7476

$ws.Range("H95").Value = 6500
$ws.Range("J95").Value = 8900
$ws.Range("L95").Value = 26700
$ws.Range("N95").Value = -30818

$ws.Range("H121").Value = 14032.77
$ws.Range("I121").Value = 18579.857
$ws.Range("J121").Value = 8727.833000000001
$ws.Range("K121").Value = 55739.571
$ws.Range("L121").Value = 26183.499
$ws.Range("M121").Value = -54429.571
$ws.Range("N121").Value = -28803.499

$ws.Range("H132").Value = 6285.5713
$ws.Range("I132").Value = 7199.8
$ws.Range("K132").Value = 64798.2
$ws.Range("M132").Value = -62268.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5722.143
$ws.Range("J70").Value = 6009.3335
$ws.Range("L70").Value = 6009.3335
$ws.Range("N70").Value = -6549.3335

$ws.Range("H73").Value = 5722.143
$ws.Range("J73").Value = 6009.3335
$ws.Range("L73").Value = 6009.3335
$ws.Range("N73").Value = -7881.3335

$ws.Range("H102").Value = 2344.2
$ws.Range("I102").Value = 2164
$ws.Range("J102").Value = 2614.5
$ws.Range("K102").Value = 2164
$ws.Range("L102").Value = 2614.5
$ws.Range("M102").Value = -542
$ws.Range("N102").Value = -5858.5

$ws.Range("H113").Value = 1098
$ws.Range("I113").Value = 1098
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1098
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1072
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 1963.8889
$ws.Range("I132").Value = 1812.8572
$ws.Range("K132").Value = 5438.571599999999
$ws.Range("M132").Value = -2908.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2727
$ws.Range("I22").Value = 2471
$ws.Range("J22").Value = 2983
$ws.Range("K22").Value = 2471
$ws.Range("L22").Value = 2983
$ws.Range("M22").Value = -2176
$ws.Range("N22").Value = -3573

$ws.Range("H26").Value = 5999
$ws.Range("I26").Value = 5999
$ws.Range("K26").Value = 5999
$ws.Range("M26").Value = -5704

$ws.Range("H27").Value = 2727
$ws.Range("I27").Value = 2471
$ws.Range("J27").Value = 2983
$ws.Range("K27").Value = 2471
$ws.Range("L27").Value = 2983
$ws.Range("M27").Value = -2364
$ws.Range("N27").Value = -3197

$ws.Range("H46").Value = 30690.861
$ws.Range("I46").Value = 64791.625
$ws.Range("J46").Value = 3410.25
$ws.Range("K46").Value = 64791.625
$ws.Range("L46").Value = 3410.25
$ws.Range("M46").Value = -64603.625
$ws.Range("N46").Value = -3786.25

$ws.Range("H68").Value = 2675.6667
$ws.Range("I68").Value = 2162.3333
$ws.Range("J68").Value = 2932.3333
$ws.Range("K68").Value = 2162.3333
$ws.Range("L68").Value = 2932.3333
$ws.Range("M68").Value = -1413.3333
$ws.Range("N68").Value = -4430.3333

$ws.Range("H71").Value = 2675.6667
$ws.Range("I71").Value = 2162.3333
$ws.Range("J71").Value = 2932.3333
$ws.Range("K71").Value = 10811.6665
$ws.Range("L71").Value = 14661.6665
$ws.Range("M71").Value = -7067.666499999999
$ws.Range("N71").Value = -22149.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H52").Value = 34999
$ws.Range("I52").Value = 29999
$ws.Range("K52").Value = 29999
$ws.Range("M52").Value = -29773

$ws.Range("H54").Value = 44999.25
$ws.Range("J54").Value = 49999.5
$ws.Range("L54").Value = 49999.5
$ws.Range("N54").Value = -51039.5

$ws.Range("H96").Value = 9647.166999999999
$ws.Range("I96").Value = 9556.200000000001
$ws.Range("J96").Value = 9712.143
$ws.Range("K96").Value = 9556.200000000001
$ws.Range("L96").Value = 9712.143
$ws.Range("M96").Value = -8183.200000000001
$ws.Range("N96").Value = -12458.143

$ws.Range("H113").Value = 536.75
$ws.Range("I113").Value = 323.75
$ws.Range("K113").Value = 971.25
$ws.Range("M113").Value = 1198.75

